$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 227
$ws.Range("I9").Value = 253.08333
$ws.Range("K9").Value = 253.08333
$ws.Range("M9").Value = -84.08332999999999
$ws.Range("H29").Value = 4522.5
$ws.Range("I29").Value = 251.5
$ws.Range("J29").Value = 5946.1665
$ws.Range("K29").Value = 754.5
$ws.Range("L29").Value = 17838.4995
$ws.Range("M29").Value = -473.5
$ws.Range("N29").Value = -18400.4995
$ws.Range("H38").Value = 1298.875
$ws.Range("J38").Value = 10000
$ws.Range("L38").Value = 30000
$ws.Range("N38").Value = -30744
$ws.Range("H40").Value = 2799.3333
$ws.Range("I40").Value = 2799.3333
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 2799.3333
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -2624.3333
$ws.Range("N40").ClearContents()
$ws.Range("H58").Value = 803.8823
$ws.Range("J58").Value = 1091.6666
$ws.Range("L58").Value = 3274.9998
$ws.Range("N58").Value = -3574.9998
$ws.Range("H64").Value = 6739.3335
$ws.Range("I64").Value = 4890.5713
$ws.Range("J64").Value = 8357
$ws.Range("K64").Value = 4890.5713
$ws.Range("L64").Value = 8357
$ws.Range("M64").Value = -4642.5713
$ws.Range("N64").Value = -8853
$ws.Range("H67").Value = 6739.3335
$ws.Range("I67").Value = 4890.5713
$ws.Range("J67").Value = 8357
$ws.Range("K67").Value = 4890.5713
$ws.Range("L67").Value = 8357
$ws.Range("M67").Value = -4032.5713
$ws.Range("N67").Value = -10073
$ws.Range("H74").Value = 7148.4
$ws.Range("I74").Value = 7246.1665
$ws.Range("K74").Value = 7246.1665
$ws.Range("M74").Value = -6310.1665
$ws.Range("H77").Value = 7148.4
$ws.Range("I77").Value = 7246.1665
$ws.Range("K77").Value = 36230.8325
$ws.Range("M77").Value = -31550.8325
$ws.Range("H86").Value = 4334.6665
$ws.Range("I86").Value = 3500
$ws.Range("K86").Value = 3500
$ws.Range("M86").Value = -2377
$ws.Range("H89").Value = 4334.6665
$ws.Range("I89").Value = 3500
$ws.Range("K89").Value = 17500
$ws.Range("M89").Value = -11884
$ws.Range("H106").Value = 3695.5334
$ws.Range("I106").Value = 2725.7693
$ws.Range("K106").Value = 2725.7693
$ws.Range("M106").Value = -2094.7693
$ws.Range("H107").Value = 3007
$ws.Range("I107").Value = 2796.2173
$ws.Range("K107").Value = 2796.2173
$ws.Range("M107").Value = -876.2172999999998
$ws.Range("H113").Value = 5266.222
$ws.Range("I113").Value = 5266.222
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 5266.222
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -2012.222
$ws.Range("N113").ClearContents()
$ws.Range("H116").Value = 10714.5
$ws.Range("I116").Value = 13191.5
$ws.Range("J116").Value = 6999
$ws.Range("K116").Value = 13191.5
$ws.Range("L116").Value = 6999
$ws.Range("M116").Value = -9749.5
$ws.Range("N116").Value = -13883
$ws.Range("H137").Value = 6452823.5
$ws.Range("I137").Value = 8334528
$ws.Range("J137").Value = 1263.8572
$ws.Range("K137").Value = 25003584
$ws.Range("L137").Value = 3791.5716
$ws.Range("M137").Value = -25001034
$ws.Range("N137").Value = -8891.571599999999
$ws.Range("H138").Value = 6373.1113
$ws.Range("I138").Value = 2195
$ws.Range("J138").Value = 6707.36
$ws.Range("K138").Value = 6585
$ws.Range("L138").Value = 20122.08
$ws.Range("M138").Value = -1445
$ws.Range("N138").Value = -30402.08
$ws.Range("H141").Value = 6928.8667
$ws.Range("I141").Value = 2869.7778
$ws.Range("J141").Value = 13017.5
$ws.Range("K141").Value = 8609.3334
$ws.Range("L141").Value = 39052.5
$ws.Range("M141").Value = -3429.3334
$ws.Range("N141").Value = -49412.5

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1368.3077
$ws.Range("J2").Value = 650
$ws.Range("L2").Value = 650
$ws.Range("N2").Value = -876
$ws.Range("H32").Value = 513854.1
$ws.Range("I32").Value = 601799.6
$ws.Range("K32").Value = 601799.6
$ws.Range("M32").Value = -601512.6
$ws.Range("H45").Value = 2108.7778
$ws.Range("I45").Value = 2247.1428
$ws.Range("K45").Value = 2247.1428
$ws.Range("M45").Value = -1870.1428
$ws.Range("H61").Value = 6690844.5
$ws.Range("I61").Value = 2527912.5
$ws.Range("K61").Value = 2527912.5
$ws.Range("M61").Value = -2527700.5
$ws.Range("H63").Value = 3401
$ws.Range("I63").Value = 3335
$ws.Range("J63").Value = 3500
$ws.Range("K63").Value = 3335
$ws.Range("L63").Value = 3500
$ws.Range("M63").Value = -2649
$ws.Range("N63").Value = -4872
$ws.Range("H66").Value = 3401
$ws.Range("I66").Value = 3335
$ws.Range("J66").Value = 3500
$ws.Range("K66").Value = 16675
$ws.Range("L66").Value = 17500
$ws.Range("M66").Value = -13243
$ws.Range("N66").Value = -24364
$ws.Range("H74").Value = 1063507.2
$ws.Range("I74").Value = 1361197.6
$ws.Range("K74").Value = 1361197.6
$ws.Range("M74").Value = -1360323.6
$ws.Range("H77").Value = 1063507.2
$ws.Range("I77").Value = 1361197.6
$ws.Range("K77").Value = 6805988
$ws.Range("M77").Value = -6801620
$ws.Range("H80").Value = 41000
$ws.Range("J80").Value = 42000
$ws.Range("L80").Value = 42000
$ws.Range("N80").Value = -43996
$ws.Range("H83").Value = 41000
$ws.Range("J83").Value = 42000
$ws.Range("L83").Value = 126000
$ws.Range("N83").Value = -135984
$ws.Range("H88").Value = 2092.923
$ws.Range("I88").Value = 1984.3334
$ws.Range("J88").Value = 2186
$ws.Range("K88").Value = 1984.3334
$ws.Range("L88").Value = 2186
$ws.Range("M88").Value = -1578.3334
$ws.Range("N88").Value = -2998
$ws.Range("H91").Value = 2092.923
$ws.Range("I91").Value = 1984.3334
$ws.Range("J91").Value = 2186
$ws.Range("K91").Value = 1984.3334
$ws.Range("L91").Value = 2186
$ws.Range("M91").Value = -580.3334
$ws.Range("N91").Value = -4994
$ws.Range("H112").Value = 157894
$ws.Range("J112").Value = 157894
$ws.Range("L112").Value = 157894
$ws.Range("N112").Value = -160848
$ws.Range("H116").Value = 1368.3077
$ws.Range("J116").Value = 650
$ws.Range("L116").Value = 650
$ws.Range("N116").Value = -5238
$ws.Range("H132").Value = 2942.1
$ws.Range("I132").Value = 1559.5333
$ws.Range("J132").Value = 7089.8
$ws.Range("K132").Value = 4678.5999
$ws.Range("L132").Value = 21269.4
$ws.Range("M132").Value = -2148.5999
$ws.Range("N132").Value = -26329.4
$ws.Range("H136").Value = 6690844.5
$ws.Range("I136").Value = 2527912.5
$ws.Range("K136").Value = 7583737.5
$ws.Range("M136").Value = -7581187.5

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1368.3077
$ws.Range("J3").Value = 650
$ws.Range("L3").Value = 650
$ws.Range("N3").Value = -878
$ws.Range("H20").Value = 41344.145
$ws.Range("I20").Value = 52397.332
$ws.Range("K20").Value = 52397.332
$ws.Range("M20").Value = -52150.332
$ws.Range("H35").Value = 37032.668
$ws.Range("J35").Value = 37032.668
$ws.Range("L35").Value = 37032.668
$ws.Range("N35").Value = -37652.668
$ws.Range("H94").Value = 2001.625
$ws.Range("I94").Value = 1760.72
$ws.Range("K94").Value = 1760.72
$ws.Range("M94").Value = -1309.72
$ws.Range("H107").Value = 3252
$ws.Range("I107").Value = 3252
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 3252
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -1332
$ws.Range("N107").ClearContents()
$ws.Range("H134").Value = 5558932.5
$ws.Range("J134").Value = 7939936.5
$ws.Range("L134").Value = 23819809.5
$ws.Range("N134").Value = -23824879.5
$ws.Range("H141").Value = 274194.75
$ws.Range("J141").Value = 274194.75
$ws.Range("L141").Value = 274194.75
$ws.Range("N141").Value = -284554.75

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H9").Value = 2397.5
$ws.Range("J9").Value = 2397.5
$ws.Range("L9").Value = 2397.5
$ws.Range("N9").Value = -2733.5
$ws.Range("H31").Value = 692401.4399999999
$ws.Range("I31").Value = 2172923.2
$ws.Range("K31").Value = 2172923.2
$ws.Range("M31").Value = -2172628.2
$ws.Range("H34").Value = 692401.4399999999
$ws.Range("I34").Value = 2172923.2
$ws.Range("K34").Value = 2172923.2
$ws.Range("M34").Value = -2172721.2
$ws.Range("H58").Value = 11834462
$ws.Range("I58").Value = 11906364
$ws.Range("J58").Value = 11498920
$ws.Range("K58").Value = 11906364
$ws.Range("L58").Value = 11498920
$ws.Range("M58").Value = -11906161
$ws.Range("N58").Value = -11499326
$ws.Range("H122").Value = 9923.888999999999
$ws.Range("I122").Value = 3142
$ws.Range("K122").Value = 9426
$ws.Range("M122").Value = -6976
$ws.Range("H132").Value = 2500.4482
$ws.Range("I132").Value = 2215.923
$ws.Range("K132").Value = 6647.768999999999
$ws.Range("M132").Value = -4117.768999999999
$ws.Range("H134").Value = 3501.725
$ws.Range("I134").Value = 1289.0952
$ws.Range("K134").Value = 3867.2856
$ws.Range("M134").Value = -1332.2856
$ws.Range("H136").Value = 11834462
$ws.Range("I136").Value = 11906364
$ws.Range("J136").Value = 11498920
$ws.Range("K136").Value = 35719092
$ws.Range("L136").Value = 34496760
$ws.Range("M136").Value = -35716542
$ws.Range("N136").Value = -34501860

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2359255.5
$ws.Range("J5").Value = 4651516
$ws.Range("L5").Value = 13954548
$ws.Range("N5").Value = -13954772
$ws.Range("H23").Value = 197.52632
$ws.Range("I23").Value = 268.57144
$ws.Range("J23").Value = 156.08333
$ws.Range("K23").Value = 805.71432
$ws.Range("L23").Value = 468.24999
$ws.Range("M23").Value = -570.71432
$ws.Range("N23").Value = -938.24999
$ws.Range("H37").Value = 95000
$ws.Range("J37").Value = 95000
$ws.Range("L37").Value = 285000
$ws.Range("N37").Value = -285224
$ws.Range("H68").Value = 6002.436
$ws.Range("I68").Value = 2000
$ws.Range("J68").Value = 6107.763
$ws.Range("K68").Value = 6000
$ws.Range("L68").Value = 18323.289
$ws.Range("M68").Value = -5189
$ws.Range("N68").Value = -19945.289
$ws.Range("H71").Value = 6002.436
$ws.Range("I71").Value = 2000
$ws.Range("J71").Value = 6107.763
$ws.Range("K71").Value = 18000
$ws.Range("L71").Value = 54969.867
$ws.Range("M71").Value = -13944
$ws.Range("N71").Value = -63081.867
$ws.Range("H98").Value = 650.2857
$ws.Range("J98").Value = 792.25
$ws.Range("L98").Value = 2376.75
$ws.Range("N98").Value = -5372.75
$ws.Range("H107").Value = 4217.7144
$ws.Range("J107").Value = 6341.278
$ws.Range("L107").Value = 19023.834
$ws.Range("N107").Value = -22863.834
$ws.Range("H114").Value = 499
$ws.Range("I114").Value = 499
$ws.Range("K114").Value = 1497
$ws.Range("M114").Value = 1757
$ws.Range("H127").Value = 5927.7856
$ws.Range("J127").Value = 5927.7856
$ws.Range("L127").Value = 17783.3568
$ws.Range("N127").Value = -27703.3568
$ws.Range("H135").Value = 2359255.5
$ws.Range("J135").Value = 4651516
$ws.Range("L135").Value = 41863644
$ws.Range("N135").Value = -41868714
$ws.Range("H137").Value = 6099.636
$ws.Range("I137").Value = 1184.8125
$ws.Range("J137").Value = 10725.353
$ws.Range("K137").Value = 3554.4375
$ws.Range("L137").Value = 32176.059
$ws.Range("M137").Value = 1545.5625
$ws.Range("N137").Value = -42376.05899999999

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 3846203.8
$ws.Range("I2").Value = 6666709
$ws.Range("K2").Value = 6666709
$ws.Range("M2").Value = -6666596
$ws.Range("H57").Value = 14739.4
$ws.Range("J57").Value = 17924.25
$ws.Range("L57").Value = 17924.25
$ws.Range("N57").Value = -19564.25
$ws.Range("H80").Value = 6931.6
$ws.Range("I80").Value = 6949.5
$ws.Range("K80").Value = 6949.5
$ws.Range("M80").Value = -5951.5
$ws.Range("H83").Value = 6931.6
$ws.Range("I83").Value = 6949.5
$ws.Range("K83").Value = 34747.5
$ws.Range("M83").Value = -29755.5
$ws.Range("H102").Value = 3029.8572
$ws.Range("I102").Value = 2677.25
$ws.Range("K102").Value = 2677.25
$ws.Range("M102").Value = -1055.25
$ws.Range("H113").Value = 1899.25
$ws.Range("I113").Value = 1953.7273
$ws.Range("J113").Value = 1300
$ws.Range("K113").Value = 1953.7273
$ws.Range("L113").Value = 1300
$ws.Range("M113").Value = 216.2727
$ws.Range("N113").Value = -5640
$ws.Range("H122").Value = 74768.71000000001
$ws.Range("I122").Value = 113918.11
$ws.Range("K122").Value = 341754.33
$ws.Range("M122").Value = -339304.33
$ws.Range("H126").Value = 2244.037
$ws.Range("I126").Value = 2257.3635
$ws.Range("K126").Value = 6772.0905
$ws.Range("M126").Value = -4302.0905
$ws.Range("H132").Value = 12093.628
$ws.Range("I132").Value = 7805.5366
$ws.Range("J132").Value = 99999.5
$ws.Range("K132").Value = 23416.6098
$ws.Range("L132").Value = 299998.5
$ws.Range("M132").Value = -20886.6098
$ws.Range("N132").Value = -305058.5
$ws.Range("H141").Value = 82085.664
$ws.Range("J141").Value = 82085.664
$ws.Range("L141").Value = 82085.664
$ws.Range("N141").Value = -92445.664

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2900
$ws.Range("I40").Value = 3298.75
$ws.Range("K40").Value = 3298.75
$ws.Range("M40").Value = -3162.75
$ws.Range("H61").Value = 6801.846
$ws.Range("I61").Value = 5179.591
$ws.Range("K61").Value = 5179.591
$ws.Range("M61").Value = -4977.591
$ws.Range("H63").Value = 49000
$ws.Range("I63").Value = 49000
$ws.Range("K63").Value = 49000
$ws.Range("M63").Value = -48251
$ws.Range("H66").Value = 49000
$ws.Range("I66").Value = 49000
$ws.Range("K66").Value = 147000
$ws.Range("M66").Value = -143256
$ws.Range("H82").Value = 5000
$ws.Range("I82").Value = 5000
$ws.Range("K82").Value = 5000
$ws.Range("M82").Value = -4639
$ws.Range("H85").Value = 5000
$ws.Range("I85").Value = 5000
$ws.Range("K85").Value = 5000
$ws.Range("M85").Value = -3752
$ws.Range("H93").Value = 10849.5
$ws.Range("I93").Value = 2000
$ws.Range("K93").Value = 2000
$ws.Range("M93").Value = -752
$ws.Range("H100").Value = 2789.3845
$ws.Range("J100").Value = 3899.8
$ws.Range("L100").Value = 3899.8
$ws.Range("N100").Value = -4981.8
$ws.Range("H110").Value = 113947
$ws.Range("J110").Value = 113947
$ws.Range("L110").Value = 113947
$ws.Range("N110").Value = -122127
$ws.Range("H113").Value = 6801.846
$ws.Range("I113").Value = 5179.591
$ws.Range("K113").Value = 5179.591
$ws.Range("M113").Value = -3009.591
$ws.Range("H122").Value = 3368.7407
$ws.Range("I122").Value = 3235.04
$ws.Range("K122").Value = 9705.119999999999
$ws.Range("M122").Value = -7255.119999999999
$ws.Range("H132").Value = 598316.25
$ws.Range("I132").Value = 903866.75
$ws.Range("J132").Value = 3296.842
$ws.Range("K132").Value = 2711600.25
$ws.Range("L132").Value = 9890.526
$ws.Range("M132").Value = -2709070.25
$ws.Range("N132").Value = -14950.526
$ws.Range("H133").Value = 85765
$ws.Range("J133").Value = 85765
$ws.Range("L133").Value = 85765
$ws.Range("N133").Value = -90825

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 20600
$ws.Range("J74").Value = 20600
$ws.Range("L74").Value = 20600
$ws.Range("N74").Value = -22472
$ws.Range("H77").Value = 20600
$ws.Range("J77").Value = 20600
$ws.Range("L77").Value = 61800
$ws.Range("N77").Value = -71160
$ws.Range("H81").Value = 5933.375
$ws.Range("I81").Value = 4417
$ws.Range("J81").Value = 8460.666999999999
$ws.Range("K81").Value = 8834
$ws.Range("L81").Value = 16921.334
$ws.Range("M81").Value = -7773
$ws.Range("N81").Value = -19043.334
$ws.Range("H84").Value = 5933.375
$ws.Range("I84").Value = 4417
$ws.Range("J84").Value = 8460.666999999999
$ws.Range("K84").Value = 44170
$ws.Range("L84").Value = 84606.67
$ws.Range("M84").Value = -38866
$ws.Range("N84").Value = -95214.67
$ws.Range("H107").Value = 1232.6666
$ws.Range("J107").Value = 799.5
$ws.Range("L107").Value = 2398.5
$ws.Range("N107").Value = -6238.5
$ws.Range("H113").Value = 2266.647
$ws.Range("I113").Value = 2144.8572
$ws.Range("K113").Value = 6434.571599999999
$ws.Range("M113").Value = -4264.571599999999
$ws.Range("H122").Value = 41862.5
$ws.Range("I122").Value = 1209.6842
$ws.Range("K122").Value = 3629.0526
$ws.Range("M122").Value = -1179.0526
$ws.Range("H132").Value = 5750340.5
$ws.Range("I132").Value = 6176034
$ws.Range("J132").Value = 3475
$ws.Range("K132").Value = 18528102
$ws.Range("L132").Value = 10425
$ws.Range("M132").Value = -18525572
$ws.Range("N132").Value = -15485
$ws.Range("H136").Value = 3072776.2
$ws.Range("I136").Value = 1787592.1
$ws.Range("K136").Value = 5362776.300000001
$ws.Range("M136").Value = -5360226.300000001
